# The workbook lists data-dictionary rows (attributeName, attributeDefinition,
# class, unit, missingValueCode, missingValueCodeExplanation) describing the
# columns of the NCP output CSV. This edit removes the two rows documenting
# "biosat" and "O2_Ar_ratio" (they are no longer part of the NCP output),
# which were rows 7 and 8 on the sheet, causing "ncp" and "k" to shift up
# from rows 9-10 to rows 7-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "biosat" (row 7) and "O2_Ar_ratio" (row 8) attribute rows.
$ws.Range("A7:G8").EntireRow.Delete() | Out-Null

# Match the author's resulting selection (rows that held "ncp"/"k" after the
# shift, i.e. what was last edited/selected).
$ws.Range("A7:XFD8").Select() | Out-Null

$wb.Save()
